$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("O2").Value = 1.23
$ws.Range("Q2").Value = 1.23
$ws.Range("S2").Value = 1.23

# Row 3
$ws.Range("N3").Value = 4.5
$ws.Range("P3").Value = 2.22
$ws.Range("S3").Value = 2.8
$ws.Range("T3").Value = 1.71
$ws.Range("U3").Value = 2.22
$ws.Range("AM3").Value = 85

# Row 4
$ws.Range("F4").Value = 4.8
$ws.Range("G4").Value = 6.2
$ws.Range("H4").Value = 1.87
$ws.Range("I4").Value = 2.08
$ws.Range("K4").Value = 3.6
$ws.Range("L4").Value = 1.48
$ws.Range("V4").Value = 1.92
$ws.Range("W4").Value = 1.2

# Row 5
$ws.Range("I5").Value = 1.6
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 3.5
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 1.89
$ws.Range("Q5").Value = 1.96
$ws.Range("R5").Value = 1.32
$ws.Range("S5").Value = 3.5
$ws.Range("T5").Value = 2.06
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 2.64
$ws.Range("X5").Value = 15

# Row 6
$ws.Range("F6").Value = 4.2
$ws.Range("H6").Value = 1.83
$ws.Range("K6").Value = 4.4
$ws.Range("P6").Value = 1.85
$ws.Range("R6").Value = 1.33
$ws.Range("U6").Value = 1.98
$ws.Range("V6").Value = 1.95
$ws.Range("Z6").Value = 1000

# Row 7
$ws.Range("F7").Value = 1.27
$ws.Range("G7").Value = 1.35
$ws.Range("I7").Value = 15
$ws.Range("K7").Value = 7.4
$ws.Range("M7").Value = 1.02
$ws.Range("O7").Value = 1.15
$ws.Range("Q7").Value = 1.44
$ws.Range("T7").Value = 1.85
$ws.Range("U7").Value = 1.94
$ws.Range("V7").Value = 1.07
$ws.Range("W7").Value = 3.85
$ws.Range("X7").Value = 40
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 13.5

# Row 8
$ws.Range("F8").Value = 1.97
$ws.Range("G8").Value = 2.18
$ws.Range("I8").Value = 4.8
$ws.Range("J8").Value = 3.55
$ws.Range("L8").Value = 1.3
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 3.95
$ws.Range("P8").Value = 2.04
$ws.Range("Q8").Value = 1.69
$ws.Range("R8").Value = 1.4
$ws.Range("S8").Value = 3
$ws.Range("U8").Value = 2.16
$ws.Range("V8").Value = 1.28
$ws.Range("W8").Value = 1.85
$ws.Range("Z8").Value = 980
$ws.Range("AD8").Value = 19.5
$ws.Range("AE8").Value = 55
$ws.Range("AN8").Value = 16

# Row 9
$ws.Range("P9").Value = 2.04
$ws.Range("Q9").Value = 1.84
$ws.Range("R9").Value = 1.41
$ws.Range("S9").Value = 3.1
$ws.Range("T9").Value = 1.71
$ws.Range("AL9").Value = 50
$ws.Range("AN9").Value = 40

# Row 10
$ws.Range("F10").Value = 3.7
$ws.Range("H10").Value = 1.79
$ws.Range("I10").Value = 2.08
$ws.Range("L10").Value = 1.27
$ws.Range("N10").Value = 4
$ws.Range("Q10").Value = 1.63
$ws.Range("R10").Value = 1.46
$ws.Range("S10").Value = 2.6
$ws.Range("T10").Value = 1.61
$ws.Range("U10").Value = 2.2
$ws.Range("V10").Value = 1.92
$ws.Range("W10").Value = 1.25
$ws.Range("AB10").Value = 23
$ws.Range("AE10").Value = 23
$ws.Range("AG10").Value = 21
$ws.Range("AH10").Value = 21

# Row 11
$ws.Range("F11").Value = 1.24
$ws.Range("G11").Value = 1.3
$ws.Range("H11").Value = 13
$ws.Range("I11").Value = 16
$ws.Range("J11").Value = 6.2
$ws.Range("K11").Value = 7.6
$ws.Range("L11").Value = 1.22
$ws.Range("N11").Value = 5.4
$ws.Range("P11").Value = 2.52
$ws.Range("R11").Value = 1.61
$ws.Range("S11").Value = 2.3
$ws.Range("T11").Value = 2.06
$ws.Range("U11").Value = 1.75
$ws.Range("V11").Value = 1.06
$ws.Range("W11").Value = 4.2
$ws.Range("Z11").Value = 170
$ws.Range("AC11").Value = 19
$ws.Range("AD11").Value = 65
$ws.Range("AE11").Value = 280
$ws.Range("AH11").Value = 40
$ws.Range("AI11").Value = 210
$ws.Range("AJ11").Value = 12
$ws.Range("AL11").Value = 48
$ws.Range("AM11").Value = 210

# Row 12
$ws.Range("G12").Value = 2.54
$ws.Range("H12").Value = 3.05
$ws.Range("L12").Value = 1.33
$ws.Range("N12").Value = 3.7
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 1.94
$ws.Range("Q12").Value = 1.89
$ws.Range("R12").Value = 1.36
$ws.Range("S12").Value = 3.3
$ws.Range("T12").Value = 1.72
$ws.Range("W12").Value = 1.65
$ws.Range("X12").Value = 18
$ws.Range("Y12").Value = 16
$ws.Range("Z12").Value = 28
$ws.Range("AA12").Value = 70
$ws.Range("AB12").Value = 13
$ws.Range("AC12").Value = 9.800000000000001
$ws.Range("AD12").Value = 17
$ws.Range("AE12").Value = 46
$ws.Range("AF12").Value = 19.5
$ws.Range("AG12").Value = 14
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 55
$ws.Range("AJ12").Value = 40
$ws.Range("AK12").Value = 32
$ws.Range("AL12").Value = 46
$ws.Range("AN12").Value = 24

# Row 13
$ws.Range("F13").Value = 1.38
$ws.Range("G13").Value = 1.63
$ws.Range("H13").Value = 7.6
$ws.Range("J13").Value = 3.55
$ws.Range("K13").Value = 7
$ws.Range("N13").Value = 1.58
$ws.Range("P13").Value = 1.58
$ws.Range("Q13").Value = 2.06
$ws.Range("S13").Value = 2.06
$ws.Range("W13").Value = 2.58

# Row 14
$ws.Range("G14").Value = 2.58
$ws.Range("I14").Value = 3.35
$ws.Range("J14").Value = 3.7
$ws.Range("M14").Value = 1.02
$ws.Range("R14").Value = 1.56
$ws.Range("S14").Value = 2.42
$ws.Range("T14").Value = 1.53
$ws.Range("V14").Value = 1.45
$ws.Range("W14").Value = 1.63

# Row 15
$ws.Range("F15").Value = 10
$ws.Range("P15").Value = 1.92

# Row 16
$ws.Range("G16").Value = 4.9
$ws.Range("K16").Value = 5
$ws.Range("P16").Value = 2.74
$ws.Range("T16").Value = 1.53
$ws.Range("Y16").Value = 18
$ws.Range("Z16").Value = 18
$ws.Range("AE16").Value = 20
$ws.Range("AH16").Value = 19.5
$ws.Range("AN16").Value = 38

# Row 17
$ws.Range("F17").Value = 2.64
$ws.Range("G17").Value = 2.96
$ws.Range("H17").Value = 2.42
$ws.Range("N17").Value = 2.38
$ws.Range("R17").Value = 1.5
$ws.Range("AC17").Value = 13

# Row 18
$ws.Range("F18").Value = 1.67
$ws.Range("H18").Value = 5.8
$ws.Range("J18").Value = 3.45
$ws.Range("T18").Value = 2.08
$ws.Range("U18").Value = 1.76
$ws.Range("AG18").Value = 12.5

# Row 20
$ws.Range("G20").Value = 1.87
$ws.Range("I20").Value = 5.5
$ws.Range("O20").Value = 1.46
$ws.Range("R20").Value = 1.25
$ws.Range("W20").Value = 2.14
$ws.Range("Y20").Value = 14.5

# Row 21
$ws.Range("G21").Value = 1.36
$ws.Range("H21").Value = 10
$ws.Range("P21").Value = 2.28
$ws.Range("R21").Value = 1.49
$ws.Range("S21").Value = 2.82
$ws.Range("T21").Value = 2.12
$ws.Range("U21").Value = 1.75
$ws.Range("AA21").Value = 460
$ws.Range("AD21").Value = 40
$ws.Range("AE21").Value = 200
$ws.Range("AH21").Value = 32
$ws.Range("AI21").Value = 170
$ws.Range("AM21").Value = 190
$ws.Range("AO21").Value = 260
